# Update the "want to go" (想去人数) counts in column F that changed
# because the underlying bilibili event data was re-scraped.
#
# The workbook has two worksheets that carry the same exhibition rows:
#   - "展览"     (Exhibitions)      -> sheet index 1
#   - "全部类型" (All types/combined) -> sheet index 4
# Both need the same updates, except row 22 ends up with a slightly
# different count between the two sheets (4238 vs 4239).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value  = 818
    $ws.Range("F3").Value  = 7
    $ws.Range("F4").Value  = 1143
    $ws.Range("F6").Value  = 12291
    $ws.Range("F7").Value  = 48
    $ws.Range("F12").Value = 897
    $ws.Range("F13").Value = 13595
    $ws.Range("F14").Value = 13762
    $ws.Range("F16").Value = 162
    $ws.Range("F19").Value = 1027
    $ws.Range("F23").Value = 209

    if ($name -eq "展览") {
        $ws.Range("F22").Value = 4238
    } else {
        $ws.Range("F22").Value = 4239
    }
}
